$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: set paragraph alignment to "right" for the Hebrew/English verse
# couplets (Gen1:2 a/b/c and Gen2:7a/b/c). For each English anchor paragraph
# we also right-align the Hebrew paragraph that immediately precedes it.
# ---------------------------------------------------------------------------
$anchors = @(
    "Now the earth was formless and empty",
    "darkness was over the surface of the deep",
    "and the Spirit of God was hovering over the waters.",
    "the Lord God formed the man from the dust of the ground",
    "breathed into his nostrils the breath of life",
    "and the man became a living being."
)

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.Trim()
    foreach ($anchor in $anchors) {
        if ($t -eq $anchor) {
            $p.Range.ParagraphFormat.Alignment = 2
            $prev = $d.Paragraphs.Item($i - 1)
            $prev.Range.ParagraphFormat.Alignment = 2
        }
    }
}
